$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 36.112135
$ws.Range("H2").Value = 108.336405
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 36.112135
$ws.Range("N2").Value = 108.336405
$ws.Range("Q2").Value = 1304.086294258225
$ws.Range("R2").Value = 11736.77664832402
